# Update crypto price ("D") and 1h-volume-change ("E") columns for rows 2-51
# to reflect the latest scrape. Values that look like plain numbers (e.g. "1.001")
# are written with a leading apostrophe so Excel keeps them as text instead of
# silently coercing them to numeric values (which would lose the trailing zero /
# thousands-dot formatting used by the source site, e.g. "27.735.36").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.735.36"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.894.40"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -1.20%  "
$ws.Range("D5").Value = "'313.29"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").Value = "'0.4837"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").Value = "'0.07340"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "'0.9180"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "'20.52"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "'0.07685"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "1.888.47"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'5.468"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "'6.599"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "'91.08"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "'0.000008794"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "27.779.38"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'14.49"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "'5.118"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "2.123.35"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "'10.78"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "'1.909"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").Value = "'153.58"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").Value = "'18.38"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "'2.131"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("D29").Value = "'115.81"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "'4.896"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "'0.08920"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'3.153"
$ws.Range("E32").Value = "  -5.55%  "
$ws.Range("D33").Value = "'1.223"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'0.7632"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").Value = "'4.647"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "'0.02042"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'2.528"
$ws.Range("E37").Value = "  -7.36%  "
$ws.Range("D38").Value = "'1.095"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").Value = "'0.05264"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "'0.5466"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").Value = "'2.987"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'6.925"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "'8.392"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "'110.09"
$ws.Range("E45").Value = "  +5.07%  "
$ws.Range("D46").Value = "'10.61"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("D50").Value = "'67.47"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "'0.06055"
$ws.Range("E51").Value = "  -0.84%  "
